# DATA GUKAR.xlsx edit
# - Add "tmt"-prep "Role" numeric codes (column J) for all data rows
# - Change E3 from "Wali Kelas" to "Bukan Wali Kelas"
# - Change F6 from "Pembina Ekstra" to "Guru Ekstra"
# - Change J6 from text "Admin" to numeric role code 1
# - Adjust a few column widths (best achievable approximation)
# - Leave selection on B7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content changes -------------------------------------------------
# (order matters: new shared-string entries are appended in first-write
# order, and the target file has "Guru Ekstra" before "Bukan Wali Kelas")
$ws.Range("F6").Value = "Guru Ekstra"
$ws.Range("E3").Value = "Bukan Wali Kelas"

# --- Column J ("Role") numeric codes for every data row (3-85) ------------
# Row 6 previously held the text "Admin" -> becomes role code 1.
# Row 24 also becomes role code 1 (newly added).
# Every other data row gets role code 4 (newly added).
for ($r = 3; $r -le 85; $r++) {
    if ($r -eq 6 -or $r -eq 24) {
        $ws.Cells.Item($r, 10).Value = 1
    } else {
        $ws.Cells.Item($r, 10).Value = 4
    }
}

# --- Column width tweaks (closest values reachable through ColumnWidth) ---
$ws.Columns.Item(6).ColumnWidth = 13.8333333333333      # F  -> ~14.71
$ws.Columns.Item(8).ColumnWidth = 9.66666666666667       # H  -> ~10.57
$ws.Columns.Item(10).ColumnWidth = 6                      # J  -> ~6.86

# --- Restore the active selection to B7 ------------------------------------
$ws.Range("B7").Select()
